$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 404,2
$arr[0,0] = "HD"
$arr[0,1] = "count"
$arr[1,0] = "HD 80586"
$arr[1,1] = 1
$arr[2,0] = "HD 87901"
$arr[2,1] = 1
$arr[3,0] = "HD 105382"
$arr[3,1] = 1
$arr[4,0] = "HD 108662"
$arr[4,1] = 1
$arr[5,0] = "HD 130841"
$arr[5,1] = 1
$arr[6,0] = "HD 11727"
$arr[6,1] = 1
$arr[7,0] = "HD 12139"
$arr[7,1] = 1
$arr[8,0] = "HD 201616"
$arr[8,1] = 1
$arr[9,0] = "HD 21743"
$arr[9,1] = 2
$arr[10,0] = "HD 21769"
$arr[10,1] = 2
$arr[11,0] = "HD 26913"
$arr[11,1] = 2
$arr[12,0] = "HD 26923"
$arr[12,1] = 2
$arr[13,0] = "HD 27778"
$arr[13,1] = 2
$arr[14,0] = "HD 28503"
$arr[14,1] = 2
$arr[15,0] = "HD 29173"
$arr[15,1] = 2
$arr[16,0] = "HD 32273"
$arr[16,1] = 2
$arr[17,0] = "HD 35149"
$arr[17,1] = 2
$arr[18,0] = "HD 42126"
$arr[18,1] = 2
$arr[19,0] = "HD 42127"
$arr[19,1] = 2
$arr[20,0] = "HD 45995"
$arr[20,1] = 2
$arr[21,0] = "HD 5788"
$arr[21,1] = 2
$arr[22,0] = "HD 5789"
$arr[22,1] = 2
$arr[23,0] = "HD 74738"
$arr[23,1] = 2
$arr[24,0] = "HD 74739"
$arr[24,1] = 2
$arr[25,0] = "HD 7439"
$arr[25,1] = 2
$arr[26,0] = "HD 82383"
$arr[26,1] = 2
$arr[27,0] = "HD 82384"
$arr[27,1] = 2
$arr[28,0] = "HD 83023"
$arr[28,1] = 2
$arr[29,0] = "HD 88849"
$arr[29,1] = 2
$arr[30,0] = "HD 105435"
$arr[30,1] = 2
$arr[31,0] = "HD 112412"
$arr[31,1] = 2
$arr[32,0] = "HD 112413"
$arr[32,1] = 2
$arr[33,0] = "HD 114447"
$arr[33,1] = 2
$arr[34,0] = "HD 129174"
$arr[34,1] = 2
$arr[35,0] = "HD 129175"
$arr[35,1] = 2
$arr[36,0] = "HD 130819"
$arr[36,1] = 2
$arr[37,0] = "HD 135101"
$arr[37,1] = 2
$arr[38,0] = "HD 139891"
$arr[38,1] = 2
$arr[39,0] = "HD 139892"
$arr[39,1] = 2
$arr[40,0] = "HD 13294"
$arr[40,1] = 2
$arr[41,0] = "HD 13295"
$arr[41,1] = 2
$arr[42,0] = "HD 166865"
$arr[42,1] = 2
$arr[43,0] = "HD 166866"
$arr[43,1] = 2
$arr[44,0] = "HD 16046"
$arr[44,1] = 2
$arr[45,0] = "HD 195093"
$arr[45,1] = 2
$arr[46,0] = "HD 195094"
$arr[46,1] = 2
$arr[47,0] = "HD 218395"
$arr[47,1] = 2
$arr[48,0] = "HD 19134"
$arr[48,1] = 2
$arr[49,0] = "HD 19135"
$arr[49,1] = 2
$arr[50,0] = "HD 24554"
$arr[50,1] = 3
$arr[51,0] = "HD 24555"
$arr[51,1] = 3
$arr[52,0] = "HD 27490"
$arr[52,1] = 3
$arr[53,0] = "HD 28271"
$arr[53,1] = 3
$arr[54,0] = "HD 30020"
$arr[54,1] = 3
$arr[55,0] = "HD 30021"
$arr[55,1] = 3
$arr[56,0] = "HD 31764"
$arr[56,1] = 3
$arr[57,0] = "HD 32039"
$arr[57,1] = 3
$arr[58,0] = "HD 32040"
$arr[58,1] = 3
$arr[59,0] = "HD 33224"
$arr[59,1] = 3
$arr[60,0] = "HD 34797"
$arr[60,1] = 3
$arr[61,0] = "HD 35736"
$arr[61,1] = 3
$arr[62,0] = "HD 36408"
$arr[62,1] = 3
$arr[63,0] = "HD 41700"
$arr[63,1] = 3
$arr[64,0] = "HD 41742"
$arr[64,1] = 3
$arr[65,0] = "HD 43017"
$arr[65,1] = 3
$arr[66,0] = "HD 4757"
$arr[66,1] = 3
$arr[67,0] = "HD 4758"
$arr[67,1] = 3
$arr[68,0] = "HD 55864"
$arr[68,1] = 3
$arr[69,0] = "HD 55865"
$arr[69,1] = 3
$arr[70,0] = "HD 57102"
$arr[70,1] = 3
$arr[71,0] = "HD 57103"
$arr[71,1] = 3
$arr[72,0] = "HD 60584"
$arr[72,1] = 3
$arr[73,0] = "HD 60585"
$arr[73,1] = 3
$arr[74,0] = "HD 6456"
$arr[74,1] = 3
$arr[75,0] = "HD 6457"
$arr[75,1] = 3
$arr[76,0] = "HD 71046"
$arr[76,1] = 3
$arr[77,0] = "HD 71150"
$arr[77,1] = 3
$arr[78,0] = "HD 71151"
$arr[78,1] = 3
$arr[79,0] = "HD 71152"
$arr[79,1] = 3
$arr[80,0] = "HD 71153"
$arr[80,1] = 3
$arr[81,0] = "HD 71176"
$arr[81,1] = 3
$arr[82,0] = "HD 71487"
$arr[82,1] = 3
$arr[83,0] = "HD 71488"
$arr[83,1] = 3
$arr[84,0] = "HD 74688"
$arr[84,1] = 3
$arr[85,0] = "HD 76369"
$arr[85,1] = 3
$arr[86,0] = "HD 76370"
$arr[86,1] = 3
$arr[87,0] = "HD 78175"
$arr[87,1] = 3
$arr[88,0] = "HD 85123"
$arr[88,1] = 3
$arr[89,0] = "HD 87344"
$arr[89,1] = 3
$arr[90,0] = "HD 91355"
$arr[90,1] = 3
$arr[91,0] = "HD 91356"
$arr[91,1] = 3
$arr[92,0] = "HD 93344"
$arr[92,1] = 3
$arr[93,0] = "HD 93359"
$arr[93,1] = 3
$arr[94,0] = "HD 94601"
$arr[94,1] = 3
$arr[95,0] = "HD 94602"
$arr[95,1] = 3
$arr[96,0] = "HD 100286"
$arr[96,1] = 3
$arr[97,0] = "HD 100287"
$arr[97,1] = 3
$arr[98,0] = "HD 105383"
$arr[98,1] = 3
$arr[99,0] = "HD 106975"
$arr[99,1] = 3
$arr[100,0] = "HD 107398"
$arr[100,1] = 3
$arr[101,0] = "HD 108651"
$arr[101,1] = 3
$arr[102,0] = "HD 108767"
$arr[102,1] = 3
$arr[103,0] = "HD 109510"
$arr[103,1] = 3
$arr[104,0] = "HD 109511"
$arr[104,1] = 3
$arr[105,0] = "HD 110379"
$arr[105,1] = 3
$arr[106,0] = "HD 110380"
$arr[106,1] = 3
$arr[107,0] = "HD 10360"
$arr[107,1] = 3
$arr[108,0] = "HD 10361"
$arr[108,1] = 3
$arr[109,0] = "HD 111720"
$arr[109,1] = 3
$arr[110,0] = "HD 112014"
$arr[110,1] = 3
$arr[111,0] = "HD 112028"
$arr[111,1] = 3
$arr[112,0] = "HD 114376"
$arr[112,1] = 3
$arr[113,0] = "HD 114846"
$arr[113,1] = 3
$arr[114,0] = "HD 116656"
$arr[114,1] = 3
$arr[115,0] = "HD 116657"
$arr[115,1] = 3
$arr[116,0] = "HD 118349"
$arr[116,1] = 3
$arr[117,0] = "HD 122408"
$arr[117,1] = 3
$arr[118,0] = "HD 124674"
$arr[118,1] = 3
$arr[119,0] = "HD 127043"
$arr[119,1] = 3
$arr[120,0] = "HD 127067"
$arr[120,1] = 3
$arr[121,0] = "HD 11503"
$arr[121,1] = 3
$arr[122,0] = "HD 11502"
$arr[122,1] = 3
$arr[123,0] = "HD 129988"
$arr[123,1] = 3
$arr[124,0] = "HD 129989"
$arr[124,1] = 3
$arr[125,0] = "HD 11749"
$arr[125,1] = 3
$arr[126,0] = "HD 134444"
$arr[126,1] = 3
$arr[127,0] = "HD 135722"
$arr[127,1] = 3
$arr[128,0] = "HD 138268"
$arr[128,1] = 3
$arr[129,0] = "HD 138488"
$arr[129,1] = 3
$arr[130,0] = "HD 139777"
$arr[130,1] = 3
$arr[131,0] = "HD 145366"
$arr[131,1] = 3
$arr[132,0] = "HD 145388"
$arr[132,1] = 3
$arr[133,0] = "HD 12533"
$arr[133,1] = 3
$arr[134,0] = "HD 12534"
$arr[134,1] = 3
$arr[135,0] = "HD 146836"
$arr[135,1] = 3
$arr[136,0] = "HD 147722"
$arr[136,1] = 3
$arr[137,0] = "HD 147723"
$arr[137,1] = 3
$arr[138,0] = "HD 150379"
$arr[138,1] = 3
$arr[139,0] = "HD 150378"
$arr[139,1] = 3
$arr[140,0] = "HD 159480"
$arr[140,1] = 3
$arr[141,0] = "HD 159541"
$arr[141,1] = 3
$arr[142,0] = "HD 159560"
$arr[142,1] = 3
$arr[143,0] = "HD 161270"
$arr[143,1] = 3
$arr[144,0] = "HD 161289"
$arr[144,1] = 3
$arr[145,0] = "HD 164668"
$arr[145,1] = 3
$arr[146,0] = "HD 164669"
$arr[146,1] = 3
$arr[147,0] = "HD 170868"
$arr[147,1] = 3
$arr[148,0] = "HD 170867"
$arr[148,1] = 3
$arr[149,0] = "HD 185644"
$arr[149,1] = 3
$arr[150,0] = "HD 187421"
$arr[150,1] = 3
$arr[151,0] = "HD 188293"
$arr[151,1] = 3
$arr[152,0] = "HD 188294"
$arr[152,1] = 3
$arr[153,0] = "HD 16232"
$arr[153,1] = 3
$arr[154,0] = "HD 16246"
$arr[154,1] = 3
$arr[155,0] = "HD 191984"
$arr[155,1] = 3
$arr[156,0] = "HD 192577"
$arr[156,1] = 3
$arr[157,0] = "HD 197963"
$arr[157,1] = 3
$arr[158,0] = "HD 197964"
$arr[158,1] = 3
$arr[159,0] = "HD 201601"
$arr[159,1] = 3
$arr[160,0] = "HD 201671"
$arr[160,1] = 3
$arr[161,0] = "HD 205811"
$arr[161,1] = 3
$arr[162,0] = "HD 209942"
$arr[162,1] = 3
$arr[163,0] = "HD 213051"
$arr[163,1] = 3
$arr[164,0] = "HD 213052"
$arr[164,1] = 3
$arr[165,0] = "HD 216172"
$arr[165,1] = 3
$arr[166,0] = "HD 219834"
$arr[166,1] = 3
$arr[167,0] = "HD 18519"
$arr[167,1] = 3
$arr[168,0] = "HD 18520"
$arr[168,1] = 3
$arr[169,0] = "HD 18622"
$arr[169,1] = 3
$arr[170,0] = "HD 18623"
$arr[170,1] = 3
$arr[171,0] = "HD 223024"
$arr[171,1] = 3
$arr[172,0] = "HD 223352"
$arr[172,1] = 3
$arr[173,0] = "HD 27710"
$arr[173,1] = 4
$arr[174,0] = "HD 31910"
$arr[174,1] = 4
$arr[175,0] = "HD 33204"
$arr[175,1] = 4
$arr[176,0] = "HD 33564"
$arr[176,1] = 4
$arr[177,0] = "HD 33959"
$arr[177,1] = 4
$arr[178,0] = "HD 34798"
$arr[178,1] = 4
$arr[179,0] = "HD 35295"
$arr[179,1] = 4
$arr[180,0] = "HD 35943"
$arr[180,1] = 4
$arr[181,0] = "HD 38393"
$arr[181,1] = 4
$arr[182,0] = "HD 44769"
$arr[182,1] = 4
$arr[183,0] = "HD 44770"
$arr[183,1] = 4
$arr[184,0] = "HD 47138"
$arr[184,1] = 4
$arr[185,0] = "HD 48250"
$arr[185,1] = 4
$arr[186,0] = "HD 48501"
$arr[186,1] = 4
$arr[187,0] = "HD 48767"
$arr[187,1] = 4
$arr[188,0] = "HD 48766"
$arr[188,1] = 4
$arr[189,0] = "HD 5156"
$arr[189,1] = 4
$arr[190,0] = "HD 53705"
$arr[190,1] = 4
$arr[191,0] = "HD 53706"
$arr[191,1] = 4
$arr[192,0] = "HD 58634"
$arr[192,1] = 4
$arr[193,0] = "HD 58635"
$arr[193,1] = 4
$arr[194,0] = "HD 62141"
$arr[194,1] = 4
$arr[195,0] = "HD 71066"
$arr[195,1] = 4
$arr[196,0] = "HD 77002"
$arr[196,1] = 4
$arr[197,0] = "HD 7344"
$arr[197,1] = 4
$arr[198,0] = "HD 7345"
$arr[198,1] = 4
$arr[199,0] = "HD 85124"
$arr[199,1] = 4
$arr[200,0] = "HD 89484"
$arr[200,1] = 4
$arr[201,0] = "HD 89485"
$arr[201,1] = 4
$arr[202,0] = "HD 89890"
$arr[202,1] = 4
$arr[203,0] = "HD 92841"
$arr[203,1] = 4
$arr[204,0] = "HD 100180"
$arr[204,1] = 4
$arr[205,0] = "HD 101177"
$arr[205,1] = 4
$arr[206,0] = "HD 103483"
$arr[206,1] = 4
$arr[207,0] = "HD 103498"
$arr[207,1] = 4
$arr[208,0] = "HD 106976"
$arr[208,1] = 4
$arr[209,0] = "HD 110317"
$arr[209,1] = 4
$arr[210,0] = "HD 110318"
$arr[210,1] = 4
$arr[211,0] = "HD 114378"
$arr[211,1] = 4
$arr[212,0] = "HD 114379"
$arr[212,1] = 4
$arr[213,0] = "HD 117200"
$arr[213,1] = 4
$arr[214,0] = "HD 117201"
$arr[214,1] = 4
$arr[215,0] = "HD 120709"
$arr[215,1] = 4
$arr[216,0] = "HD 120710"
$arr[216,1] = 4
$arr[217,0] = "HD 126128"
$arr[217,1] = 4
$arr[218,0] = "HD 126129"
$arr[218,1] = 4
$arr[219,0] = "HD 126367"
$arr[219,1] = 4
$arr[220,0] = "HD 131977"
$arr[220,1] = 4
$arr[221,0] = "HD 133408"
$arr[221,1] = 4
$arr[222,0] = "HD 134443"
$arr[222,1] = 4
$arr[223,0] = "HD 11973"
$arr[223,1] = 4
$arr[224,0] = "HD 137107"
$arr[224,1] = 4
$arr[225,0] = "HD 137108"
$arr[225,1] = 4
$arr[226,0] = "HD 137391"
$arr[226,1] = 4
$arr[227,0] = "HD 137392"
$arr[227,1] = 4
$arr[228,0] = "HD 139460"
$arr[228,1] = 4
$arr[229,0] = "HD 139461"
$arr[229,1] = 4
$arr[230,0] = "HD 12446"
$arr[230,1] = 4
$arr[231,0] = "HD 12447"
$arr[231,1] = 4
$arr[232,0] = "HD 144069"
$arr[232,1] = 4
$arr[233,0] = "HD 144070"
$arr[233,1] = 4
$arr[234,0] = "HD 144217"
$arr[234,1] = 4
$arr[235,0] = "HD 144218"
$arr[235,1] = 4
$arr[236,0] = "HD 145001"
$arr[236,1] = 4
$arr[237,0] = "HD 145000"
$arr[237,1] = 4
$arr[238,0] = "HD 146362"
$arr[238,1] = 4
$arr[239,0] = "HD 154905"
$arr[239,1] = 4
$arr[240,0] = "HD 154906"
$arr[240,1] = 4
$arr[241,0] = "HD 156014"
$arr[241,1] = 4
$arr[242,0] = "HD 156015"
$arr[242,1] = 4
$arr[243,0] = "HD 166045"
$arr[243,1] = 4
$arr[244,0] = "HD 166046"
$arr[244,1] = 4
$arr[245,0] = "HD 177463"
$arr[245,1] = 4
$arr[246,0] = "HD 179957"
$arr[246,1] = 4
$arr[247,0] = "HD 179958"
$arr[247,1] = 4
$arr[248,0] = "HD 187420"
$arr[248,1] = 4
$arr[249,0] = "HD 190147"
$arr[249,1] = 4
$arr[250,0] = "HD 191570"
$arr[250,1] = 4
$arr[251,0] = "HD 192514"
$arr[251,1] = 4
$arr[252,0] = "HD 198160"
$arr[252,1] = 4
$arr[253,0] = "HD 198161"
$arr[253,1] = 4
$arr[254,0] = "HD 206826"
$arr[254,1] = 4
$arr[255,0] = "HD 206827"
$arr[255,1] = 4
$arr[256,0] = "HD 209278"
$arr[256,1] = 4
$arr[257,0] = "HD 18537"
$arr[257,1] = 4
$arr[258,0] = "HD 18538"
$arr[258,1] = 4
$arr[259,0] = "HD 224635"
$arr[259,1] = 4
$arr[260,0] = "HD 224636"
$arr[260,1] = 4
$arr[261,0] = "HD 225009"
$arr[261,1] = 4
$arr[262,0] = "HD 24071"
$arr[262,1] = 5
$arr[263,0] = "HD 24072"
$arr[263,1] = 5
$arr[264,0] = "HD 2884"
$arr[264,1] = 5
$arr[265,0] = "HD 2885"
$arr[265,1] = 5
$arr[266,0] = "HD 28446"
$arr[266,1] = 5
$arr[267,0] = "HD 35162"
$arr[267,1] = 5
$arr[268,0] = "HD 38392"
$arr[268,1] = 5
$arr[269,0] = "HD 57852"
$arr[269,1] = 5
$arr[270,0] = "HD 57853"
$arr[270,1] = 5
$arr[271,0] = "HD 62153"
$arr[271,1] = 5
$arr[272,0] = "HD 62154"
$arr[272,1] = 5
$arr[273,0] = "HD 62863"
$arr[273,1] = 5
$arr[274,0] = "HD 62864"
$arr[274,1] = 5
$arr[275,0] = "HD 74560"
$arr[275,1] = 5
$arr[276,0] = "HD 82780"
$arr[276,1] = 5
$arr[277,0] = "HD 97855"
$arr[277,1] = 5
$arr[278,0] = "HD 99491"
$arr[278,1] = 5
$arr[279,0] = "HD 112092"
$arr[279,1] = 5
$arr[280,0] = "HD 112091"
$arr[280,1] = 5
$arr[281,0] = "HD 115810"
$arr[281,1] = 5
$arr[282,0] = "HD 125161"
$arr[282,1] = 5
$arr[283,0] = "HD 129246"
$arr[283,1] = 5
$arr[284,0] = "HD 129247"
$arr[284,1] = 5
$arr[285,0] = "HD 129926"
$arr[285,1] = 5
$arr[286,0] = "HD 142630"
$arr[286,1] = 5
$arr[287,0] = "HD 146361"
$arr[287,1] = 5
$arr[288,0] = "HD 150117"
$arr[288,1] = 5
$arr[289,0] = "HD 150118"
$arr[289,1] = 5
$arr[290,0] = "HD 152909"
$arr[290,1] = 5
$arr[291,0] = "HD 157778"
$arr[291,1] = 5
$arr[292,0] = "HD 157779"
$arr[292,1] = 5
$arr[293,0] = "HD 164764"
$arr[293,1] = 5
$arr[294,0] = "HD 164765"
$arr[294,1] = 5
$arr[295,0] = "HD 173582"
$arr[295,1] = 5
$arr[296,0] = "HD 173583"
$arr[296,1] = 5
$arr[297,0] = "HD 175638"
$arr[297,1] = 5
$arr[298,0] = "HD 175639"
$arr[298,1] = 5
$arr[299,0] = "HD 199766"
$arr[299,1] = 5
$arr[300,0] = "HD 208095"
$arr[300,1] = 5
$arr[301,0] = "HD 213306"
$arr[301,1] = 5
$arr[302,0] = "HD 214599"
$arr[302,1] = 5
$arr[303,0] = "HD 219449"
$arr[303,1] = 5
$arr[304,0] = "HD 28255"
$arr[304,1] = 6
$arr[305,0] = "HD 37646"
$arr[305,1] = 6
$arr[306,0] = "HD 60178"
$arr[306,1] = 6
$arr[307,0] = "HD 60179"
$arr[307,1] = 6
$arr[308,0] = "HD 68257"
$arr[308,1] = 6
$arr[309,0] = "HD 68256"
$arr[309,1] = 6
$arr[310,0] = "HD 68255"
$arr[310,1] = 6
$arr[311,0] = "HD 72946"
$arr[311,1] = 6
$arr[312,0] = "HD 75086"
$arr[312,1] = 6
$arr[313,0] = "HD 98230"
$arr[313,1] = 6
$arr[314,0] = "HD 98231"
$arr[314,1] = 6
$arr[315,0] = "HD 108903"
$arr[315,1] = 6
$arr[316,0] = "HD 124675"
$arr[316,1] = 6
$arr[317,0] = "HD 138917"
$arr[317,1] = 6
$arr[318,0] = "HD 138918"
$arr[318,1] = 6
$arr[319,0] = "HD 142629"
$arr[319,1] = 6
$arr[320,0] = "HD 143118"
$arr[320,1] = 6
$arr[321,0] = "HD 150100"
$arr[321,1] = 6
$arr[322,0] = "HD 156349"
$arr[322,1] = 6
$arr[323,0] = "HD 156350"
$arr[323,1] = 6
$arr[324,0] = "HD 162003"
$arr[324,1] = 6
$arr[325,0] = "HD 162004"
$arr[325,1] = 6
$arr[326,0] = "HD 165189"
$arr[326,1] = 6
$arr[327,0] = "HD 165190"
$arr[327,1] = 6
$arr[328,0] = "HD 173648"
$arr[328,1] = 6
$arr[329,0] = "HD 173649"
$arr[329,1] = 6
$arr[330,0] = "HD 174638"
$arr[330,1] = 6
$arr[331,0] = "HD 176269"
$arr[331,1] = 6
$arr[332,0] = "HD 176270"
$arr[332,1] = 6
$arr[333,0] = "HD 187013"
$arr[333,1] = 6
$arr[334,0] = "HD 209790"
$arr[334,1] = 6
$arr[335,0] = "HD 45725"
$arr[335,1] = 7
$arr[336,0] = "HD 45726"
$arr[336,1] = 7
$arr[337,0] = "HD 45727"
$arr[337,1] = 7
$arr[338,0] = "HD 59067"
$arr[338,1] = 7
$arr[339,0] = "HD 59499"
$arr[339,1] = 7
$arr[340,0] = "HD 59500"
$arr[340,1] = 7
$arr[341,0] = "HD 6479"
$arr[341,1] = 7
$arr[342,0] = "HD 6480"
$arr[342,1] = 7
$arr[343,0] = "HD 66005"
$arr[343,1] = 7
$arr[344,0] = "HD 72945"
$arr[344,1] = 7
$arr[345,0] = "HD 74535"
$arr[345,1] = 7
$arr[346,0] = "HD 116072"
$arr[346,1] = 7
$arr[347,0] = "HD 116087"
$arr[347,1] = 7
$arr[348,0] = "HD 140483"
$arr[348,1] = 7
$arr[349,0] = "HD 140484"
$arr[349,1] = 7
$arr[350,0] = "HD 173087"
$arr[350,1] = 7
$arr[351,0] = "HD 173607"
$arr[351,1] = 7
$arr[352,0] = "HD 173608"
$arr[352,1] = 7
$arr[353,0] = "HD 186901"
$arr[353,1] = 7
$arr[354,0] = "HD 201091"
$arr[354,1] = 7
$arr[355,0] = "HD 201092"
$arr[355,1] = 7
$arr[356,0] = "HD 212697"
$arr[356,1] = 7
$arr[357,0] = "HD 212698"
$arr[357,1] = 7
$arr[358,0] = "HD 66006"
$arr[358,1] = 8
$arr[359,0] = "HD 7788"
$arr[359,1] = 8
$arr[360,0] = "HD 135734"
$arr[360,1] = 8
$arr[361,0] = "HD 36485"
$arr[361,1] = 9
$arr[362,0] = "HD 147553"
$arr[362,1] = 9
$arr[363,0] = "HD 155885"
$arr[363,1] = 9
$arr[364,0] = "HD 155886"
$arr[364,1] = 9
$arr[365,0] = "HD 183914"
$arr[365,1] = 10
$arr[366,0] = "HD 36486"
$arr[366,1] = 11
$arr[367,0] = "HD 108248"
$arr[367,1] = 11
$arr[368,0] = "HD 108925"
$arr[368,1] = 11
$arr[369,0] = "HD 108249"
$arr[369,1] = 12
$arr[370,0] = "HD 147933"
$arr[370,1] = 12
$arr[371,0] = "HD 147934"
$arr[371,1] = 12
$arr[372,0] = "HD 145501"
$arr[372,1] = 13
$arr[373,0] = "HD 145502"
$arr[373,1] = 13
$arr[374,0] = "HD 186408"
$arr[374,1] = 14
$arr[375,0] = "HD 186427"
$arr[375,1] = 14
$arr[376,0] = "HD 108250"
$arr[376,1] = 15
$arr[377,0] = "HD 183912"
$arr[377,1] = 15
$arr[378,0] = "HD 214168"
$arr[378,1] = 15
$arr[379,0] = "HD 37742"
$arr[379,1] = 16
$arr[380,0] = "HD 37743"
$arr[380,1] = 16
$arr[381,0] = "HD 61555"
$arr[381,1] = 16
$arr[382,0] = "HD 61556"
$arr[382,1] = 18
$arr[383,0] = "HD 128620"
$arr[383,1] = 18
$arr[384,0] = "HD 128621"
$arr[384,1] = 18
$arr[385,0] = "HD 144668"
$arr[385,1] = 20
$arr[386,0] = "HD 144667"
$arr[386,1] = 21
$arr[387,0] = "HD 68243"
$arr[387,1] = 22
$arr[388,0] = "HD 68273"
$arr[388,1] = 25
$arr[389,0] = "HD 36959"
$arr[389,1] = 29
$arr[390,0] = "HD 36960"
$arr[390,1] = 29
$arr[391,0] = "HD 25638"
$arr[391,1] = 39
$arr[392,0] = "HD 163755"
$arr[392,1] = 42
$arr[393,0] = "HD 163756"
$arr[393,1] = 43
$arr[394,0] = "HD 36861"
$arr[394,1] = 62
$arr[395,0] = "HD 36862"
$arr[395,1] = 63
$arr[396,0] = "HD 37479"
$arr[396,1] = 71
$arr[397,0] = "HD 37468"
$arr[397,1] = 72
$arr[398,0] = "HD 150136"
$arr[398,1] = 100
$arr[399,0] = "HD 37041"
$arr[399,1] = 578
$arr[400,0] = "HD 37023"
$arr[400,1] = 1906
$arr[401,0] = "HD 37022"
$arr[401,1] = 1927
$arr[402,0] = "HD 37020"
$arr[402,1] = 1937
$arr[403,0] = "HD 37021"
$arr[403,1] = 1938

$ws.Range("A1:B404").Value = $arr